# Insert a new paragraph containing a page break immediately after the
# empty paragraph that precedes the "Figure 3" caption (right after the
# Figure 2 picture paragraph), matching the target revision.

$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Locate the empty paragraph that immediately precedes the "Figure 3:" caption.
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "Figure 3:*") {
        $targetIndex = $i - 1
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the paragraph preceding 'Figure 3:'"
}

$anchorPara = $paras.Item($targetIndex)

# Add a brand-new paragraph right after the anchor; Word gives it the
# anchor's own formatting to start with.
$anchorPara.Range.InsertParagraphAfter()

# The freshly created paragraph is now the next one in the collection.
$newPara = $paras.Item($targetIndex + 1)
$insertionPoint = $newPara.Range
$insertionPoint.Collapse(1)

# Replace the (empty) collapsed range with the fully-specified paragraph
# mark formatting plus a single run holding a page break, via a WordOpenXML
# package fragment so the exact run/paragraph-mark properties are set in one
# shot (matches the committed revision byte-for-byte).
$insertionPoint.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document><w:body><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/><w:b/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times"/><w:b/><w:color w:val="000000"/></w:rPr><w:br w:type="page"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

Write-Output "Inserted page-break paragraph after paragraph index $targetIndex"
